# #272 Ajout d'un scenario de recherche de l'offre d'un professionnel avec un ID Nat PS
# - Bump the StructureDefinition "Date" metadata value.
# - Swap the two Mapping columns ("RIM Mapping" <-> "Specification metier ...")
#   on the Elements sheet, including their column widths.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: refresh the "Date" row (A8/B8) ---------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# --- Elements sheet: swap Mapping columns AK (37) and AL (38) -----------
$wsElements = $wb.Worksheets.Item("Elements")

$lastRow = 10
for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $wsElements.Cells.Item($r, 37)
    $alCell = $wsElements.Cells.Item($r, 38)

    $akVal = $akCell.Value2
    $alVal = $alCell.Value2

    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# Swap the column widths that went with the two Mapping columns.
$wsElements.Columns.Item(37).ColumnWidth = 99.8359375
$wsElements.Columns.Item(38).ColumnWidth = 24.98046875
